# Actualización automática 2025-08-18 15:50:08
#
# The August sale figure for client "MUÑOZ LOZA ROMMEL SEBASTIAN" (advisor
# HIDALGO HIDALGO PEDRO GUSTAVO) under PORCELANATO increased by 12.86, which
# ripples into the per-sheet monthly / category totals and the compliance
# summary sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" --------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M14").Value = 12.86
$wsGrupo.Range("M23").Value = "3 de 21"

# --- Sheet "VENTA MENSUAL" ------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F14").Value = 142.46
$wsMensual.Range("F23").Value = 2909.72

# --- Sheet "CUMPLIMIENTO MENSUAL" -----------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Excel's ColumnWidth property adds a fixed 5/6-character padding on top of
# the raw stored <col width> value, so back it out to land on exactly 13.
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13 - 5/6

$wsCumplimiento.Range("D16").Value = 2668.76
$wsCumplimiento.Range("E16").Value = 36107.71
$wsCumplimiento.Range("F16").Value = 0.06882421220910516

$wsCumplimiento.Range("D19").Value = 2909.72
$wsCumplimiento.Range("E19").Value = 56478.50762291768
$wsCumplimiento.Range("F19").Value = 0.04899489539366469
